$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 'Switzerland, Government of'
$ws.Cells.Item(2, 2).Value = 12207540
$ws.Cells.Item(3, 1).Value = 'Japan, Government of'
$ws.Cells.Item(3, 2).Value = 5047309
$ws.Cells.Item(4, 1).Value = 'United States of America, Government of'
$ws.Cells.Item(4, 2).Value = 133652732
$ws.Cells.Item(5, 1).Value = 'Walmart Foundation'
$ws.Cells.Item(5, 2).Value = 500000
$ws.Cells.Item(6, 1).Value = 'Saint Lucia, Government of'
$ws.Cells.Item(6, 2).Value = 0
$ws.Cells.Item(7, 1).Value = 'Barbados, Government of'
$ws.Cells.Item(7, 2).Value = 0
$ws.Cells.Item(8, 1).Value = 'Direct Relief'
$ws.Cells.Item(8, 2).Value = 0
$ws.Cells.Item(9, 1).Value = 'European Commission'
$ws.Cells.Item(9, 2).Value = 5107592
$ws.Cells.Item(10, 1).Value = 'Jamaica, Government of'
$ws.Cells.Item(10, 2).Value = 0
$ws.Cells.Item(11, 1).Value = 'Turks and Caicos Islands'
$ws.Cells.Item(11, 2).Value = 0
$ws.Cells.Item(12, 1).Value = 'Canada, Government of'
$ws.Cells.Item(12, 2).Value = 30482351
$ws.Cells.Item(13, 1).Value = 'Mexico, Government of'
$ws.Cells.Item(13, 2).Value = 25000
$ws.Cells.Item(14, 1).Value = 'Private (individuals & organizations)'
$ws.Cells.Item(14, 2).Value = 1245089
$ws.Cells.Item(15, 1).Value = 'FedEx'
$ws.Cells.Item(15, 2).Value = 0
$ws.Cells.Item(16, 1).Value = 'PepsiCo Foundation'
$ws.Cells.Item(16, 2).Value = 0
$ws.Cells.Item(17, 1).Value = 'Western Union Foundation'
$ws.Cells.Item(17, 2).Value = 25000
$ws.Cells.Item(18, 1).Value = 'Agility'
$ws.Cells.Item(18, 2).Value = 0
$ws.Cells.Item(19, 1).Value = 'Qatar Charity'
$ws.Cells.Item(19, 2).Value = 56962
$ws.Cells.Item(20, 1).Value = 'Liechtenstein, Government of'
$ws.Cells.Item(20, 2).Value = 252274
$ws.Cells.Item(21, 1).Value = 'European Commission''s Humanitarian Aid and Civil Protection Department'
$ws.Cells.Item(21, 2).Value = 31867844
$ws.Cells.Item(22, 1).Value = 'Central Emergency Response Fund'
$ws.Cells.Item(22, 2).Value = 20318394
$ws.Cells.Item(23, 1).Value = 'United Kingdom, Government of'
$ws.Cells.Item(23, 2).Value = 10445383
$ws.Cells.Item(24, 1).Value = 'Ireland, Government of'
$ws.Cells.Item(24, 2).Value = 2324258
$ws.Cells.Item(25, 1).Value = 'Accion Contra el Hambre - Spain'
$ws.Cells.Item(25, 2).Value = 928221
$ws.Cells.Item(26, 1).Value = 'Germany, Government of'
$ws.Cells.Item(26, 2).Value = 27975629
$ws.Cells.Item(27, 1).Value = 'World Vision Taiwan'
$ws.Cells.Item(27, 2).Value = 284562
$ws.Cells.Item(28, 1).Value = 'World Vision Australia'
$ws.Cells.Item(28, 2).Value = 346196
$ws.Cells.Item(29, 1).Value = 'Denmark, Government of'
$ws.Cells.Item(29, 2).Value = 375400
$ws.Cells.Item(30, 1).Value = 'Argentina, Government of'
$ws.Cells.Item(30, 2).Value = 1305563
$ws.Cells.Item(31, 1).Value = 'France, Government of'
$ws.Cells.Item(31, 2).Value = 8952678
$ws.Cells.Item(32, 1).Value = 'Brazil, Government of'
$ws.Cells.Item(32, 2).Value = 323866
$ws.Cells.Item(33, 1).Value = 'Estonia, Government of'
$ws.Cells.Item(33, 2).Value = 163501
$ws.Cells.Item(34, 1).Value = 'Qatar, Government of'
$ws.Cells.Item(34, 2).Value = 800000
$ws.Cells.Item(35, 1).Value = 'Luxembourg, Government of'
$ws.Cells.Item(35, 2).Value = 1237968
$ws.Cells.Item(36, 1).Value = 'Sweden, Government of'
$ws.Cells.Item(36, 2).Value = 9386422
$ws.Cells.Item(37, 1).Value = 'National Bank of Canada'
$ws.Cells.Item(37, 2).Value = 50000
$ws.Cells.Item(38, 1).Value = 'UNICEF National Committee/France'
$ws.Cells.Item(38, 2).Value = 1548873
$ws.Cells.Item(39, 1).Value = 'Italy, Government of'
$ws.Cells.Item(39, 2).Value = 1121076
$ws.Cells.Item(40, 1).Value = 'Australia, Government of'
$ws.Cells.Item(40, 2).Value = 1900391
$ws.Cells.Item(41, 1).Value = 'New Zealand, Government of'
$ws.Cells.Item(41, 2).Value = 354610
$ws.Cells.Item(42, 1).Value = 'UNICEF National Committee/Canada'
$ws.Cells.Item(42, 2).Value = 172015
$ws.Cells.Item(43, 1).Value = 'UNICEF National Committee/Netherlands'
$ws.Cells.Item(43, 2).Value = 54348
$ws.Cells.Item(44, 1).Value = 'UNICEF National Committee/Denmark'
$ws.Cells.Item(44, 2).Value = 718094
$ws.Cells.Item(45, 1).Value = 'UNICEF National Committee/Germany'
$ws.Cells.Item(45, 2).Value = 824358
$ws.Cells.Item(46, 1).Value = 'UNICEF National Committee/Japan'
$ws.Cells.Item(46, 2).Value = 15063
$ws.Cells.Item(47, 1).Value = 'UNICEF National Committee/Spain'
$ws.Cells.Item(47, 2).Value = 2810410
$ws.Cells.Item(48, 1).Value = 'US Fund for UNICEF'
$ws.Cells.Item(48, 2).Value = 2465158
$ws.Cells.Item(49, 1).Value = 'Norway, Government of'
$ws.Cells.Item(49, 2).Value = 845002
$ws.Cells.Item(50, 1).Value = 'Spain, Government of'
$ws.Cells.Item(50, 2).Value = 1352660
$ws.Cells.Item(51, 1).Value = 'UNICEF National Committee/Korea (Republic of)'
$ws.Cells.Item(51, 2).Value = 51362
$ws.Cells.Item(52, 1).Value = 'UNICEF National Committee/Norway'
$ws.Cells.Item(52, 2).Value = 28676
$ws.Cells.Item(53, 1).Value = 'UNICEF National Committee/United Kingdom'
$ws.Cells.Item(53, 2).Value = 31095
$ws.Cells.Item(54, 1).Value = 'UNICEF National Committee/Czech Republic'
$ws.Cells.Item(54, 2).Value = 10949
$ws.Cells.Item(55, 1).Value = 'UNICEF National Committee/Andorra'
$ws.Cells.Item(55, 2).Value = 11792
$ws.Cells.Item(56, 1).Value = 'Korea, Republic of, Government of'
$ws.Cells.Item(56, 2).Value = 200000
$ws.Cells.Item(57, 1).Value = 'Lithuania, Government of'
$ws.Cells.Item(57, 2).Value = 21858
$ws.Cells.Item(58, 1).Value = 'Mastercard Foundation'
$ws.Cells.Item(58, 2).Value = 100000
$ws.Cells.Item(59, 1).Value = 'Slovenia, Government of'
$ws.Cells.Item(59, 2).Value = 54645
$ws.Cells.Item(60, 1).Value = ''
$ws.Cells.Item(60, 2).Value = 1230593
$ws.Cells.Item(61, 1).Value = 'World Food Programme'
$ws.Cells.Item(61, 2).Value = 2600804
$ws.Cells.Item(62, 1).Value = 'Finland, Government of'
$ws.Cells.Item(62, 2).Value = 68370
$ws.Cells.Item(63, 1).Value = 'United Arab Emirates, Government of'
$ws.Cells.Item(63, 2).Value = 400000
$ws.Cells.Item(64, 1).Value = 'Swiss Solidarity'
$ws.Cells.Item(64, 2).Value = 4024426
$ws.Cells.Item(65, 1).Value = 'Microsoft'
$ws.Cells.Item(65, 2).Value = 193388
$ws.Cells.Item(66, 1).Value = 'UPS Foundation'
$ws.Cells.Item(66, 2).Value = 100000
$ws.Cells.Item(67, 1).Value = 'United Nations Population Fund'
$ws.Cells.Item(67, 2).Value = 255880
$ws.Cells.Item(68, 1).Value = 'Friends of UNFPA'
$ws.Cells.Item(68, 2).Value = 1210733
$ws.Cells.Item(69, 1).Value = 'Baby Box Company'
$ws.Cells.Item(69, 2).Value = 112500
$ws.Cells.Item(70, 1).Value = 'UN Foundation'
$ws.Cells.Item(70, 2).Value = 99997
$ws.Cells.Item(71, 1).Value = 'Chile, Government of'
$ws.Cells.Item(71, 2).Value = 50000
$ws.Cells.Item(72, 1).Value = 'Monaco, Government of'
$ws.Cells.Item(72, 2).Value = 60109
$ws.Cells.Item(73, 1).Value = 'Singapore, Government of'
$ws.Cells.Item(73, 2).Value = 50000
$ws.Cells.Item(74, 1).Value = 'Romania, Government of'
$ws.Cells.Item(74, 2).Value = 55188
$ws.Cells.Item(75, 1).Value = 'Malta, Government of'
$ws.Cells.Item(75, 2).Value = 27174
$ws.Cells.Item(76, 1).Value = 'Thailand, Government of'
$ws.Cells.Item(76, 2).Value = 100000
$ws.Cells.Item(77, 1).Value = 'Peru, Government of'
$ws.Cells.Item(77, 2).Value = 10319
$ws.Cells.Item(78, 1).Value = 'Office for the Coordination of Humanitarian Affairs'
$ws.Cells.Item(78, 2).Value = 66696
$ws.Cells.Item(79, 1).Value = 'International Organization for Migration'
$ws.Cells.Item(79, 2).Value = 689792
$ws.Cells.Item(80, 1).Value = 'Botswana, Government of'
$ws.Cells.Item(80, 2).Value = 50000
$ws.Cells.Item(81, 1).Value = 'Food & Agriculture Organization of the United Nations'
$ws.Cells.Item(81, 2).Value = 500000
$ws.Cells.Item(82, 1).Value = 'Belgium, Government of'
$ws.Cells.Item(82, 2).Value = 3778688
$ws.Cells.Item(83, 1).Value = 'United Nations Development Programme'
$ws.Cells.Item(83, 2).Value = 500000
$ws.Cells.Item(84, 1).Value = 'World Jewish Relief'
$ws.Cells.Item(84, 2).Value = 32510
$ws.Cells.Item(85, 1).Value = 'Medicor Foundation'
$ws.Cells.Item(85, 2).Value = 240000
$ws.Cells.Item(86, 1).Value = 'Starbucks Foundation'
$ws.Cells.Item(86, 2).Value = 25000
$ws.Cells.Item(87, 1).Value = 'Nuestros Pequeños Hermanos'
$ws.Cells.Item(87, 2).Value = 9893113
